$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (Property/Value table) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row;
# turn it into "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second duplicate "Contact" row - remove it entirely,
# shifting all subsequent rows up by one.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" (element definitions table) ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (root Extension element) Short/Definition columns (K/L)
$elements.Range("K2").Value = "Quality Measures"
$elements.Range("L2").Value = "Quality measures related to the FFV intitiative"
